$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) from row 2 to row 27: change 45208 -> 45212
for ($r = 2; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45208) {
        $cell.Value2 = 45212
    }
}
